$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The title-row merge is getting wider (C1:F1 -> C1:G1) to make room for the
# new "Mã danh mục" column, so unmerge first and re-merge at the end.
$ws.Range("C1:F1").UnMerge()

# Row 2 - header: insert the new "Mã danh mục" column after "Tên sản phẩm"
# and shift the remaining headers one column to the right.
$ws.Cells.Item(2, 3).Value = "Mã danh mục"
$ws.Cells.Item(2, 4).Value = "Tồn đầu kỳ"
$ws.Cells.Item(2, 5).Value = "Số lượng nhập"
$ws.Cells.Item(2, 6).Value = "Số lượng xuất"
$ws.Cells.Item(2, 7).Value = "Tồn cuối kỳ"

# Rows 3-8 - shift existing quantity data one column right and fill the new
# "Mã danh mục" column with 1.
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 6
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 4

$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0

$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0

$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0

$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0

$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0

# Row 9 (new) - "may giat 1" product.
$ws.Cells.Item(9, 1).Value = 34
$ws.Cells.Item(9, 2).Value = "may giat 1"
$ws.Cells.Item(9, 3).Value = 3
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0

# New column G needs the same rendered width (10) as the other data
# columns. The engine's ColumnWidth setter adds the usual Excel
# padding (~0.8333 chars) on top of what we ask for, so back it out here.
$ws.Columns("G").ColumnWidth = 9.166666666666666

# Re-merge the title row across the new, wider range.
$ws.Range("C1:G1").Merge()
